# Cel-Cxcr4 LR-pair sheet: refresh NATMI ligand/receptor metrics with the
# re-run (TPM-based) numbers for rows 2-17 (columns E-J ligand-side,
# M-T receptor-side / edge-weight derived-specificity columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 0.834641
$ws.Range("H2").Value2 = 2.503923
$ws.Range("I2").Value2 = 0.4389248184909427
$ws.Range("J2").Value2 = 0.4389248184909426
$ws.Range("M2").Value2 = 30.52246933333333
$ws.Range("N2").Value2 = 91.567408
$ws.Range("O2").Value2 = 0.1058764512547768
$ws.Range("P2").Value2 = 0.1058764512547769
$ws.Range("Q2").Value2 = 25.47530432684266
$ws.Range("R2").Value2 = 229.277738941584
$ws.Range("S2").Value2 = 0.04647180214946808
$ws.Range("T2").Value2 = 0.04647180214946807
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 0.834641
$ws.Range("H3").Value2 = 2.503923
$ws.Range("I3").Value2 = 0.4389248184909427
$ws.Range("J3").Value2 = 0.4389248184909426
$ws.Range("O3").Value2 = 0.001067503492562006
$ws.Range("P3").Value2 = 0.001067503492562006
$ws.Range("Q3").Value2 = 0.2568557599039999
$ws.Range("R3").Value2 = 2.311701839136
$ws.Range("S3").Value2 = 0.0004685537767112257
$ws.Range("T3").Value2 = 0.0004685537767112257
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 0.834641
$ws.Range("H4").Value2 = 2.503923
$ws.Range("I4").Value2 = 0.4389248184909427
$ws.Range("J4").Value2 = 0.4389248184909426
$ws.Range("M4").Value2 = 47.57542166666667
$ws.Range("N4").Value2 = 142.726265
$ws.Range("O4").Value2 = 0.1650297935598315
$ws.Range("P4").Value2 = 0.1650297935598315
$ws.Range("Q4").Value2 = 39.70839751528833
$ws.Range("R4").Value2 = 357.375577637595
$ws.Range("S4").Value2 = 0.07243567218384681
$ws.Range("T4").Value2 = 0.07243567218384678
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 0.834641
$ws.Range("H5").Value2 = 2.503923
$ws.Range("I5").Value2 = 0.4389248184909427
$ws.Range("J5").Value2 = 0.4389248184909426
$ws.Range("M5").Value2 = 209.8781993333333
$ws.Range("N5").Value2 = 629.634598
$ws.Range("O5").Value2 = 0.7280262516928295
$ws.Range("P5").Value2 = 0.7280262516928296
$ws.Range("Q5").Value2 = 175.1729501697727
$ws.Range("R5").Value2 = 1576.556551527954
$ws.Range("S5").Value2 = 0.3195487903809166
$ws.Range("T5").Value2 = 0.3195487903809165
$ws.Range("G6").Value2 = 0.3724893333333333
$ws.Range("I6").Value2 = 0.1958863907034828
$ws.Range("J6").Value2 = 0.1958863907034828
$ws.Range("M6").Value2 = 30.52246933333333
$ws.Range("N6").Value2 = 91.567408
$ws.Range("O6").Value2 = 0.1058764512547768
$ws.Range("P6").Value2 = 0.1058764512547769
$ws.Range("Q6").Value2 = 11.36929425366044
$ws.Range("R6").Value2 = 102.323648282944
$ws.Range("S6").Value2 = 0.02073975589679147
$ws.Range("T6").Value2 = 0.02073975589679147
$ws.Range("G7").Value2 = 0.3724893333333333
$ws.Range("I7").Value2 = 0.1958863907034828
$ws.Range("J7").Value2 = 0.1958863907034828
$ws.Range("O7").Value2 = 0.001067503492562006
$ws.Range("P7").Value2 = 0.001067503492562006
$ws.Range("S7").Value2 = 0.0002091094062213334
$ws.Range("T7").Value2 = 0.0002091094062213335
$ws.Range("G8").Value2 = 0.3724893333333333
$ws.Range("I8").Value2 = 0.1958863907034828
$ws.Range("J8").Value2 = 0.1958863907034828
$ws.Range("M8").Value2 = 47.57542166666667
$ws.Range("N8").Value2 = 142.726265
$ws.Range("O8").Value2 = 0.1650297935598315
$ws.Range("P8").Value2 = 0.1650297935598315
$ws.Range("Q8").Value2 = 17.72133709966889
$ws.Range("R8").Value2 = 159.49203389702
$ws.Range("S8").Value2 = 0.03232709061897627
$ws.Range("T8").Value2 = 0.03232709061897626
$ws.Range("G9").Value2 = 0.3724893333333333
$ws.Range("I9").Value2 = 0.1958863907034828
$ws.Range("J9").Value2 = 0.1958863907034828
$ws.Range("M9").Value2 = 209.8781993333333
$ws.Range("N9").Value2 = 629.634598
$ws.Range("O9").Value2 = 0.7280262516928295
$ws.Range("P9").Value2 = 0.7280262516928296
$ws.Range("Q9").Value2 = 78.17739055087377
$ws.Range("R9").Value2 = 703.596514957864
$ws.Range("S9").Value2 = 0.1426104347814937
$ws.Range("T9").Value2 = 0.1426104347814937
$ws.Range("G10").Value2 = 0.4120993333333334
$ws.Range("H10").Value2 = 1.236298
$ws.Range("I10").Value2 = 0.2167166782887156
$ws.Range("J10").Value2 = 0.2167166782887155
$ws.Range("M10").Value2 = 30.52246933333333
$ws.Range("N10").Value2 = 91.567408
$ws.Range("O10").Value2 = 0.1058764512547768
$ws.Range("P10").Value2 = 0.1058764512547769
$ws.Range("Q10").Value2 = 12.57828926395378
$ws.Range("R10").Value2 = 113.204603375584
$ws.Range("S10").Value2 = 0.02294519282493235
$ws.Range("T10").Value2 = 0.02294519282493235
$ws.Range("G11").Value2 = 0.4120993333333334
$ws.Range("H11").Value2 = 1.236298
$ws.Range("I11").Value2 = 0.2167166782887156
$ws.Range("J11").Value2 = 0.2167166782887155
$ws.Range("O11").Value2 = 0.001067503492562006
$ws.Range("P11").Value2 = 0.001067503492562006
$ws.Range("Q11").Value2 = 0.1268210972373333
$ws.Range("R11").Value2 = 1.141389875136
$ws.Range("S11").Value2 = 0.0002313458109696404
$ws.Range("T11").Value2 = 0.0002313458109696404
$ws.Range("G12").Value2 = 0.4120993333333334
$ws.Range("H12").Value2 = 1.236298
$ws.Range("I12").Value2 = 0.2167166782887156
$ws.Range("J12").Value2 = 0.2167166782887155
$ws.Range("M12").Value2 = 47.57542166666667
$ws.Range("N12").Value2 = 142.726265
$ws.Range("O12").Value2 = 0.1650297935598315
$ws.Range("P12").Value2 = 0.1650297935598315
$ws.Range("Q12").Value2 = 19.60579955188556
$ws.Range("R12").Value2 = 176.45219596697
$ws.Range("S12").Value2 = 0.03576470867895915
$ws.Range("T12").Value2 = 0.03576470867895915
$ws.Range("G13").Value2 = 0.4120993333333334
$ws.Range("H13").Value2 = 1.236298
$ws.Range("I13").Value2 = 0.2167166782887156
$ws.Range("J13").Value2 = 0.2167166782887155
$ws.Range("M13").Value2 = 209.8781993333333
$ws.Range("N13").Value2 = 629.634598
$ws.Range("O13").Value2 = 0.7280262516928295
$ws.Range("P13").Value2 = 0.7280262516928296
$ws.Range("Q13").Value2 = 86.49066602646711
$ws.Range("R13").Value2 = 778.415994238204
$ws.Range("S13").Value2 = 0.1577754309738544
$ws.Range("T13").Value2 = 0.1577754309738544
$ws.Range("G14").Value2 = 0.2823283333333333
$ws.Range("H14").Value2 = 0.8469849999999999
$ws.Range("I14").Value2 = 0.148472112516859
$ws.Range("J14").Value2 = 0.1484721125168589
$ws.Range("M14").Value2 = 30.52246933333333
$ws.Range("N14").Value2 = 91.567408
$ws.Range("O14").Value2 = 0.1058764512547768
$ws.Range("P14").Value2 = 0.1058764512547769
$ws.Range("Q14").Value2 = 8.617357896097776
$ws.Range("R14").Value2 = 77.55622106487999
$ws.Range("S14").Value2 = 0.01571970038358496
$ws.Range("T14").Value2 = 0.01571970038358496
$ws.Range("G15").Value2 = 0.2823283333333333
$ws.Range("H15").Value2 = 0.8469849999999999
$ws.Range("I15").Value2 = 0.148472112516859
$ws.Range("J15").Value2 = 0.1484721125168589
$ws.Range("O15").Value2 = 0.001067503492562006
$ws.Range("P15").Value2 = 0.001067503492562006
$ws.Range("Q15").Value2 = 0.08688485061333331
$ws.Range("R15").Value2 = 0.7819636555199998
$ws.Range("S15").Value2 = 0.000158494498659806
$ws.Range("T15").Value2 = 0.000158494498659806
$ws.Range("G16").Value2 = 0.2823283333333333
$ws.Range("H16").Value2 = 0.8469849999999999
$ws.Range("I16").Value2 = 0.148472112516859
$ws.Range("J16").Value2 = 0.1484721125168589
$ws.Range("M16").Value2 = 47.57542166666667
$ws.Range("N16").Value2 = 142.726265
$ws.Range("O16").Value2 = 0.1650297935598315
$ws.Range("P16").Value2 = 0.1650297935598315
$ws.Range("Q16").Value2 = 13.43188950678056
$ws.Range("R16").Value2 = 120.887005561025
$ws.Range("S16").Value2 = 0.02450232207804931
$ws.Range("T16").Value2 = 0.02450232207804931
$ws.Range("G17").Value2 = 0.2823283333333333
$ws.Range("H17").Value2 = 0.8469849999999999
$ws.Range("I17").Value2 = 0.148472112516859
$ws.Range("J17").Value2 = 0.1484721125168589
$ws.Range("M17").Value2 = 209.8781993333333
$ws.Range("N17").Value2 = 629.634598
$ws.Range("O17").Value2 = 0.7280262516928295
$ws.Range("P17").Value2 = 0.7280262516928296
$ws.Range("Q17").Value2 = 59.2545622207811
$ws.Range("R17").Value2 = 533.2910599870299
$ws.Range("S17").Value2 = 0.1080915955565649
$ws.Range("T17").Value2 = 0.1080915955565649
